$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D601").Value = 45194
$ws.Range("I601").Value = "Primera"
$ws.Range("J601").Value = 5000
$ws.Range("K601").Value = 500
$ws.Range("L601").Value = 500
$ws.Range("M601").Value = 500
$ws.Range("P601").Value = 100
$ws.Range("D602").Value = 44532
$ws.Range("I602").Value = "Primera"
$ws.Range("J602").Value = 5000
$ws.Range("K602").Value = 500
$ws.Range("L602").Value = 500
$ws.Range("M602").Value = 500
$ws.Range("P602").Value = 100
$ws.Range("D603").Value = 44701
$ws.Range("I603").Value = "Primera"
$ws.Range("J603").Value = 4000
$ws.Range("K603").Value = 600
$ws.Range("L603").Value = 600
$ws.Range("M603").Value = 600
$ws.Range("P603").Value = 120
$ws.Range("D604").Value = 44999
$ws.Range("I604").Value = "Primera"
$ws.Range("J604").Value = 5000
$ws.Range("K604").Value = 600
$ws.Range("L604").Value = 600
$ws.Range("M604").Value = 600
$ws.Range("P604").Value = 120
$ws.Range("D605").Value = 44543
$ws.Range("I605").Value = "Primera"
$ws.Range("J605").Value = 5000
$ws.Range("K605").Value = 500
$ws.Range("L605").Value = 500
$ws.Range("M605").Value = 500
$ws.Range("P605").Value = 100
$ws.Range("D606").Value = 45075
$ws.Range("I606").Value = "Primera"
$ws.Range("J606").Value = 3000
$ws.Range("K606").Value = 600
$ws.Range("L606").Value = 600
$ws.Range("M606").Value = 600
$ws.Range("P606").Value = 120
$ws.Range("D607").Value = 45075
$ws.Range("I607").Value = "Segunda"
$ws.Range("J607").Value = 3000
$ws.Range("K607").Value = 500
$ws.Range("L607").Value = 500
$ws.Range("M607").Value = 500
$ws.Range("P607").Value = 100
$ws.Range("D608").Value = 44572
$ws.Range("I608").Value = "Primera"
$ws.Range("J608").Value = 4000
$ws.Range("K608").Value = 500
$ws.Range("L608").Value = 500
$ws.Range("M608").Value = 500
$ws.Range("P608").Value = 100
$ws.Range("D609").Value = 45063
$ws.Range("I609").Value = "Primera"
$ws.Range("J609").Value = 3000
$ws.Range("K609").Value = 600
$ws.Range("L609").Value = 600
$ws.Range("M609").Value = 600
$ws.Range("P609").Value = 120
$ws.Range("D610").Value = 45063
$ws.Range("I610").Value = "Segunda"
$ws.Range("J610").Value = 2000
$ws.Range("K610").Value = 500
$ws.Range("L610").Value = 500
$ws.Range("M610").Value = 500
$ws.Range("P610").Value = 100
$ws.Range("D611").Value = 44201
$ws.Range("I611").Value = "Primera"
$ws.Range("J611").Value = 3000
$ws.Range("K611").Value = 500
$ws.Range("L611").Value = 500
$ws.Range("M611").Value = 500
$ws.Range("P611").Value = 100
$ws.Range("D612").Value = 44776
$ws.Range("I612").Value = "Primera"
$ws.Range("J612").Value = 4000
$ws.Range("K612").Value = 750
$ws.Range("L612").Value = 750
$ws.Range("M612").Value = 750
$ws.Range("P612").Value = 150
$ws.Range("D613").Value = 44315
$ws.Range("I613").Value = "Primera"
$ws.Range("J613").Value = 3000
$ws.Range("K613").Value = 600
$ws.Range("L613").Value = 600
$ws.Range("M613").Value = 600
$ws.Range("P613").Value = 120
$ws.Range("D614").Value = 44432
$ws.Range("I614").Value = "Primera"
$ws.Range("J614").Value = 3000
$ws.Range("K614").Value = 600
$ws.Range("L614").Value = 600
$ws.Range("M614").Value = 600
$ws.Range("P614").Value = 120
$ws.Range("D615").Value = 44295
$ws.Range("I615").Value = "Primera"
$ws.Range("J615").Value = 3000
$ws.Range("K615").Value = 600
$ws.Range("L615").Value = 600
$ws.Range("M615").Value = 600
$ws.Range("P615").Value = 120
$ws.Range("D616").Value = 44466
$ws.Range("I616").Value = "Primera"
$ws.Range("J616").Value = 4000
$ws.Range("K616").Value = 700
$ws.Range("L616").Value = 700
$ws.Range("M616").Value = 700
$ws.Range("P616").Value = 140
$ws.Range("D617").Value = 44980
$ws.Range("I617").Value = "Primera"
$ws.Range("J617").Value = 5000
$ws.Range("K617").Value = 600
$ws.Range("L617").Value = 600
$ws.Range("M617").Value = 600
$ws.Range("P617").Value = 120
$ws.Range("D618").Value = 44517
$ws.Range("I618").Value = "Primera"
$ws.Range("J618").Value = 5000
$ws.Range("K618").Value = 550
$ws.Range("L618").Value = 550
$ws.Range("M618").Value = 550
$ws.Range("P618").Value = 110
$ws.Range("D619").Value = 44614
$ws.Range("I619").Value = "Primera"
$ws.Range("J619").Value = 2000
$ws.Range("K619").Value = 800
$ws.Range("L619").Value = 800
$ws.Range("M619").Value = 800
$ws.Range("P619").Value = 160
$ws.Range("D620").Value = 45027
$ws.Range("I620").Value = "Primera"
$ws.Range("J620").Value = 4000
$ws.Range("K620").Value = 600
$ws.Range("L620").Value = 600
$ws.Range("M620").Value = 600
$ws.Range("P620").Value = 120
$ws.Range("D621").Value = 45027
$ws.Range("I621").Value = "Segunda"
$ws.Range("J621").Value = 2000
$ws.Range("K621").Value = 500
$ws.Range("L621").Value = 500
$ws.Range("M621").Value = 500
$ws.Range("P621").Value = 100
$ws.Range("D622").Value = 44708
$ws.Range("I622").Value = "Primera"
$ws.Range("J622").Value = 4000
$ws.Range("K622").Value = 650
$ws.Range("L622").Value = 650
$ws.Range("M622").Value = 650
$ws.Range("P622").Value = 130
$ws.Range("D623").Value = 44449
$ws.Range("I623").Value = "Primera"
$ws.Range("J623").Value = 3000
$ws.Range("K623").Value = 600
$ws.Range("L623").Value = 600
$ws.Range("M623").Value = 600
$ws.Range("P623").Value = 120
$ws.Range("D624").Value = 44223
$ws.Range("I624").Value = "Primera"
$ws.Range("J624").Value = 3000
$ws.Range("K624").Value = 600
$ws.Range("L624").Value = 600
$ws.Range("M624").Value = 600
$ws.Range("P624").Value = 120
$ws.Range("D625").Value = 44781
$ws.Range("I625").Value = "Primera"
$ws.Range("J625").Value = 3000
$ws.Range("K625").Value = 700
$ws.Range("L625").Value = 700
$ws.Range("M625").Value = 700
$ws.Range("P625").Value = 140
$ws.Range("D626").Value = 44706
$ws.Range("I626").Value = "Primera"
$ws.Range("J626").Value = 5000
$ws.Range("K626").Value = 600
$ws.Range("L626").Value = 600
$ws.Range("M626").Value = 600
$ws.Range("P626").Value = 120
$ws.Range("D627").Value = 44251
$ws.Range("I627").Value = "Primera"
$ws.Range("J627").Value = 2000
$ws.Range("K627").Value = 600
$ws.Range("L627").Value = 600
$ws.Range("M627").Value = 600
$ws.Range("P627").Value = 120
$ws.Range("D628").Value = 44676
$ws.Range("I628").Value = "Primera"
$ws.Range("J628").Value = 5000
$ws.Range("K628").Value = 600
$ws.Range("L628").Value = 600
$ws.Range("M628").Value = 600
$ws.Range("P628").Value = 120
$ws.Range("D629").Value = 44600
$ws.Range("I629").Value = "Segunda"
$ws.Range("J629").Value = 3000
$ws.Range("K629").Value = 700
$ws.Range("L629").Value = 700
$ws.Range("M629").Value = 700
$ws.Range("P629").Value = 140
$ws.Range("D630").Value = 44484
$ws.Range("I630").Value = "Primera"
$ws.Range("J630").Value = 4000
$ws.Range("K630").Value = 700
$ws.Range("L630").Value = 700
$ws.Range("M630").Value = 700
$ws.Range("P630").Value = 140
$ws.Range("D631").Value = 44609
$ws.Range("I631").Value = "Segunda"
$ws.Range("J631").Value = 3000
$ws.Range("K631").Value = 700
$ws.Range("L631").Value = 700
$ws.Range("M631").Value = 700
$ws.Range("P631").Value = 140
$ws.Range("D632").Value = 44518
$ws.Range("I632").Value = "Primera"
$ws.Range("J632").Value = 6000
$ws.Range("K632").Value = 500
$ws.Range("L632").Value = 500
$ws.Range("M632").Value = 500
$ws.Range("P632").Value = 100
$ws.Range("D633").Value = 44756
$ws.Range("I633").Value = "Primera"
$ws.Range("J633").Value = 5000
$ws.Range("K633").Value = 700
$ws.Range("L633").Value = 700
$ws.Range("M633").Value = 700
$ws.Range("P633").Value = 140
$ws.Range("D634").Value = 44467
$ws.Range("I634").Value = "Primera"
$ws.Range("J634").Value = 4000
$ws.Range("K634").Value = 650
$ws.Range("L634").Value = 650
$ws.Range("M634").Value = 650
$ws.Range("P634").Value = 130
$ws.Range("D635").Value = 45076
$ws.Range("I635").Value = "Primera"
$ws.Range("J635").Value = 3000
$ws.Range("K635").Value = 600
$ws.Range("L635").Value = 600
$ws.Range("M635").Value = 600
$ws.Range("P635").Value = 120
$ws.Range("D636").Value = 45076
$ws.Range("I636").Value = "Segunda"
$ws.Range("J636").Value = 2000
$ws.Range("K636").Value = 500
$ws.Range("L636").Value = 500
$ws.Range("M636").Value = 500
$ws.Range("P636").Value = 100
$ws.Range("D637").Value = 44327
$ws.Range("I637").Value = "Primera"
$ws.Range("J637").Value = 5000
$ws.Range("K637").Value = 500
$ws.Range("L637").Value = 500
$ws.Range("M637").Value = 500
$ws.Range("P637").Value = 100
$ws.Range("D638").Value = 44664
$ws.Range("I638").Value = "Primera"
$ws.Range("J638").Value = 6000
$ws.Range("K638").Value = 600
$ws.Range("L638").Value = 600
$ws.Range("M638").Value = 600
$ws.Range("P638").Value = 120
$ws.Range("D639").Value = 44460
$ws.Range("I639").Value = "Primera"
$ws.Range("J639").Value = 3000
$ws.Range("K639").Value = 650
$ws.Range("L639").Value = 650
$ws.Range("M639").Value = 650
$ws.Range("P639").Value = 130
$ws.Range("D640").Value = 45012
$ws.Range("I640").Value = "Primera"
$ws.Range("J640").Value = 5000
$ws.Range("K640").Value = 650
$ws.Range("L640").Value = 650
$ws.Range("M640").Value = 650
$ws.Range("P640").Value = 130
$ws.Range("D641").Value = 44364
$ws.Range("I641").Value = "Primera"
$ws.Range("J641").Value = 5000
$ws.Range("K641").Value = 500
$ws.Range("L641").Value = 500
$ws.Range("M641").Value = 500
$ws.Range("P641").Value = 100
$ws.Range("D642").Value = 44222
$ws.Range("I642").Value = "Primera"
$ws.Range("J642").Value = 3000
$ws.Range("K642").Value = 600
$ws.Range("L642").Value = 600
$ws.Range("M642").Value = 600
$ws.Range("P642").Value = 120
$ws.Range("D643").Value = 44711
$ws.Range("I643").Value = "Primera"
$ws.Range("J643").Value = 4000
$ws.Range("K643").Value = 650
$ws.Range("L643").Value = 650
$ws.Range("M643").Value = 650
$ws.Range("P643").Value = 130
$ws.Range("D644").Value = 44279
$ws.Range("I644").Value = "Primera"
$ws.Range("J644").Value = 3000
$ws.Range("K644").Value = 500
$ws.Range("L644").Value = 500
$ws.Range("M644").Value = 500
$ws.Range("P644").Value = 100
$ws.Range("D645").Value = 44715
$ws.Range("I645").Value = "Primera"
$ws.Range("J645").Value = 3000
$ws.Range("K645").Value = 650
$ws.Range("L645").Value = 650
$ws.Range("M645").Value = 650
$ws.Range("P645").Value = 130
$ws.Range("D646").Value = 44965
$ws.Range("I646").Value = "Primera"
$ws.Range("J646").Value = 5000
$ws.Range("K646").Value = 500
$ws.Range("L646").Value = 500
$ws.Range("M646").Value = 500
$ws.Range("P646").Value = 100
$ws.Range("D647").Value = 45069
$ws.Range("I647").Value = "Primera"
$ws.Range("J647").Value = 3000
$ws.Range("K647").Value = 600
$ws.Range("L647").Value = 600
$ws.Range("M647").Value = 600
$ws.Range("P647").Value = 120
$ws.Range("D648").Value = 45069
$ws.Range("I648").Value = "Segunda"
$ws.Range("J648").Value = 2000
$ws.Range("K648").Value = 500
$ws.Range("L648").Value = 500
$ws.Range("M648").Value = 500
$ws.Range("P648").Value = 100
$ws.Range("D649").Value = 44540
$ws.Range("I649").Value = "Primera"
$ws.Range("J649").Value = 5000
$ws.Range("K649").Value = 500
$ws.Range("L649").Value = 500
$ws.Range("M649").Value = 500
$ws.Range("P649").Value = 100
$ws.Range("D650").Value = 44362
$ws.Range("I650").Value = "Primera"
$ws.Range("J650").Value = 6000
$ws.Range("K650").Value = 500
$ws.Range("L650").Value = 500
$ws.Range("M650").Value = 500
$ws.Range("P650").Value = 100
$ws.Range("D651").Value = 44826
$ws.Range("I651").Value = "Segunda"
$ws.Range("J651").Value = 3000
$ws.Range("K651").Value = 800
$ws.Range("L651").Value = 800
$ws.Range("M651").Value = 800
$ws.Range("P651").Value = 160
$ws.Range("D652").Value = 44845
$ws.Range("I652").Value = "Segunda"
$ws.Range("J652").Value = 4000
$ws.Range("K652").Value = 800
$ws.Range("L652").Value = 800
$ws.Range("M652").Value = 800
$ws.Range("P652").Value = 160
$ws.Range("D653").Value = 44567
$ws.Range("I653").Value = "Primera"
$ws.Range("J653").Value = 5000
$ws.Range("K653").Value = 500
$ws.Range("L653").Value = 500
$ws.Range("M653").Value = 500
$ws.Range("P653").Value = 100
$ws.Range("D654").Value = 44525
$ws.Range("I654").Value = "Primera"
$ws.Range("J654").Value = 5000
$ws.Range("K654").Value = 500
$ws.Range("L654").Value = 500
$ws.Range("M654").Value = 500
$ws.Range("P654").Value = 100
$ws.Range("D655").Value = 44656
$ws.Range("I655").Value = "Primera"
$ws.Range("J655").Value = 5000
$ws.Range("K655").Value = 600
$ws.Range("L655").Value = 600
$ws.Range("M655").Value = 600
$ws.Range("P655").Value = 120
$ws.Range("D656").Value = 45149
$ws.Range("I656").Value = "Primera"
$ws.Range("J656").Value = 5000
$ws.Range("K656").Value = 500
$ws.Range("L656").Value = 500
$ws.Range("M656").Value = 500
$ws.Range("P656").Value = 100
$ws.Range("D657").Value = 44382
$ws.Range("I657").Value = "Primera"
$ws.Range("J657").Value = 6000
$ws.Range("K657").Value = 600
$ws.Range("L657").Value = 600
$ws.Range("M657").Value = 600
$ws.Range("P657").Value = 120
$ws.Range("D658").Value = 44557
$ws.Range("I658").Value = "Primera"
$ws.Range("J658").Value = 5000
$ws.Range("K658").Value = 550
$ws.Range("L658").Value = 550
$ws.Range("M658").Value = 550
$ws.Range("P658").Value = 110
$ws.Range("D659").Value = 45175
$ws.Range("I659").Value = "Primera"
$ws.Range("J659").Value = 3000
$ws.Range("K659").Value = 500
$ws.Range("L659").Value = 500
$ws.Range("M659").Value = 500
$ws.Range("P659").Value = 100
$ws.Range("D660").Value = 45173
$ws.Range("I660").Value = "Primera"
$ws.Range("J660").Value = 6000
$ws.Range("K660").Value = 450
$ws.Range("L660").Value = 500
$ws.Range("M660").Value = 475
$ws.Range("P660").Value = 95
$ws.Range("D661").Value = 44200
$ws.Range("I661").Value = "Primera"
$ws.Range("J661").Value = 3000
$ws.Range("K661").Value = 500
$ws.Range("L661").Value = 500
$ws.Range("M661").Value = 500
$ws.Range("P661").Value = 100
$ws.Range("D662").Value = 44991
$ws.Range("I662").Value = "Primera"
$ws.Range("J662").Value = 5000
$ws.Range("K662").Value = 600
$ws.Range("L662").Value = 600
$ws.Range("M662").Value = 600
$ws.Range("P662").Value = 120
$ws.Range("D663").Value = 45191
$ws.Range("I663").Value = "Primera"
$ws.Range("J663").Value = 5000
$ws.Range("K663").Value = 500
$ws.Range("L663").Value = 500
$ws.Range("M663").Value = 500
$ws.Range("P663").Value = 100

# New row 663 - static metadata columns copied from row 662
$ws.Range("A663").Value = 5
$ws.Range("B663").Value = "Macroferia Regional de Talca"
$ws.Range("C663").Value = "Maule"
$ws.Range("E663").Value = 7
$ws.Range("F663").Value = 100114014
$ws.Range("G663").Value = "Betarraga"
$ws.Range("H663").Value = "Sin especificar"
$ws.Range("N663").Value = "$/paquete 5 unidades"
$ws.Range("O663").Value = "Región del Maule"
$ws.Range("Q663").Value = 5
$ws.Range("R663").Value = "Hortaliza"
$ws.Range("D663").NumberFormat = "YYYY-MM-DD HH:MM:SS"